$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.296.34'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '2.491.26'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = "'321.04"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.49%  '
$ws.Range('D6').Value = "'108.41"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.06%  '
$ws.Range('D7').Value = "'0.522"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('E9').Value = '  -1.11%  '
$ws.Range('D10').Value = "'39.47"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.79%  '
$ws.Range('D11').Value = "'0.0810"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.55%  '
$ws.Range('E12').Value = '  +0.22%  '
$ws.Range('D13').Value = "'18.35"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.07%  '
$ws.Range('E14').Value = '  -1.31%  '
$ws.Range('D15').Value = '2.883.14'
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('D16').Value = '2.503.73'
$ws.Range('E16').Value = '  +0.90%  '
$ws.Range('D17').Value = "'0.846"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.08%  '
$ws.Range('D18').Value = '47.212.69'
$ws.Range('E18').Value = '  +0.39%  '
$ws.Range('D19').Value = "'13.12"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.88%  '
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('E22').Value = '  +11.93%  '
$ws.Range('D23').Value = "'70.35"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.82%  '
$ws.Range('D24').Value = "'245.04"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.54%  '
$ws.Range('E25').Value = '  +0.73%  '
$ws.Range('D26').Value = "'1.00"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').Value = "'25.70"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.85%  '
$ws.Range('D28').Value = "'2.28"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.62%  '
$ws.Range('E29').Value = '  -1.62%  '
$ws.Range('D30').Value = "'34.75"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.78%  '
$ws.Range('E31').Value = '  +1.65%  '
$ws.Range('D32').Value = "'49.80"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.28%  '
$ws.Range('D33').Value = "'20.61"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.31%  '
$ws.Range('D34').Value = "'5.35"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.47%  '
$ws.Range('D35').Value = "'0.0784"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').Value = "'4.72"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.25%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').Value = "'1.97"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.84%  '
$ws.Range('E39').Value = '  -1.62%  '
$ws.Range('D40').Value = "'22.93"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.05%  '
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('E42').Value = '  +0.51%  '
$ws.Range('D43').Value = "'116.70"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.41%  '
$ws.Range('D44').Value = "'0.0296"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.68%  '
$ws.Range('D45').Value = '1.997.41'
$ws.Range('E45').Value = '  +2.21%  '
$ws.Range('E46').Value = '  +1.84%  '
$ws.Range('D47').Value = "'1.99"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.20%  '
$ws.Range('E48').Value = '  +0.40%  '
$ws.Range('E49').Value = '  -0.63%  '
$ws.Range('E50').Value = '  -5.09%  '
$ws.Range('D51').Value = "'56.41"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.16%  '
